$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Append a duplicate of the most recent week's block (rows 109-112,
#        dated 44392) as new rows 113-116, preserving all original values.
$ws.Range("A109:T112").Copy()
$ws.Range("A113").PasteSpecial(-4104)
# PasteSpecial creates a fresh date-formatted style for column D instead of
# reusing the existing one (style index 2) - put it back so the new cells
# match the rest of the date column.
$ws.Range("D113:D116").NumberFormat = $ws.Range("D109").NumberFormat

# --- 2) Shift the remaining weekly blocks down: every block keeps its own
#        price figures, only its reporting date ("D" column) moves to the
#        next block's date - except the newest block (rows 101-104), which
#        gets a brand-new date (44476) together with the prices that used
#        to belong to the now-duplicated last block.

# Row 101 (Especial)
$ws.Range("D101").Value = 44476
$ws.Range("N101").Value = 18000
$ws.Range("O101").Value = 19000
$ws.Range("P101").Value = 18500
$ws.Range("S101").Value = 1850

# Row 102 (Primera)
$ws.Range("D102").Value = 44476
$ws.Range("N102").Value = 18000
$ws.Range("O102").Value = 19000
$ws.Range("P102").Value = 18500
$ws.Range("S102").Value = 1542

# Row 103 (Segunda)
$ws.Range("D103").Value = 44476
$ws.Range("N103").Value = 18000
$ws.Range("O103").Value = 19000
$ws.Range("P103").Value = 18500
$ws.Range("S103").Value = 1321

# Row 104 (Tercera)
$ws.Range("D104").Value = 44476
$ws.Range("N104").Value = 18000
$ws.Range("O104").Value = 19000
$ws.Range("P104").Value = 18500
$ws.Range("S104").Value = 1156

# Row 105 (Especial) - date only
$ws.Range("D105").Value = 44306

# Row 106 (Primera) - date only
$ws.Range("D106").Value = 44306

# Row 107 (Segunda) - date only
$ws.Range("D107").Value = 44306

# Row 108 (Tercera) - date only
$ws.Range("D108").Value = 44306

# Row 109 (Especial)
$ws.Range("D109").Value = 44357
$ws.Range("N109").Value = 14000
$ws.Range("O109").Value = 15000
$ws.Range("P109").Value = 14500
$ws.Range("S109").Value = 1450

# Row 110 (Primera)
$ws.Range("D110").Value = 44357
$ws.Range("M110").Value = 250
$ws.Range("N110").Value = 14000
$ws.Range("O110").Value = 15000
$ws.Range("P110").Value = 14500
$ws.Range("S110").Value = 1208

# Row 111 (Segunda)
$ws.Range("D111").Value = 44357
$ws.Range("N111").Value = 14000
$ws.Range("O111").Value = 15000
$ws.Range("P111").Value = 14500
$ws.Range("S111").Value = 1036

# Row 112 (Tercera)
$ws.Range("D112").Value = 44357
$ws.Range("N112").Value = 14000
$ws.Range("O112").Value = 15000
$ws.Range("P112").Value = 14500
$ws.Range("S112").Value = 906
